$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Housekeeping" table (header row + 3 data rows) that used to sit on
# rows 45-48 is moved up by one row, to rows 44-47 (freeing up row 48).
$ws.Range("A45:S45").Copy($ws.Range("A44"))
$ws.Range("A46:S46").Copy($ws.Range("A45"))
$ws.Range("A47:S47").Copy($ws.Range("A46"))
$ws.Range("A48:S48").Copy($ws.Range("A47"))

# Column A only holds a value on the (copied) header row now; clear the
# stale leftovers on the data rows that the copy above doesn't blank out.
$ws.Range("A45:A47").ClearContents()

# Drop the now-empty trailing row so the sheet's used range shrinks back
# down to row 47.
$ws.Rows("48:48").Delete()

# Range.Copy above carries values/number-formats/styles but not the
# underlying boolean formulas - restore those explicitly on the three
# data rows so I:S keep their TRUE()/FALSE() formulas.
$cols = @("I","J","K","L","M","N","O","P","Q","R","S")

$row45 = @($true,$false,$false,$false,$true,$true,$true,$true,$true,$true,$false)
$row46 = @($true,$true,$false,$false,$true,$true,$true,$true,$true,$true,$false)
$row47 = @($true,$false,$false,$false,$true,$false,$false,$true,$true,$true,$true)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $fn = "FALSE()"
    if ($row45[$i]) { $fn = "TRUE()" }
    $ws.Range($cols[$i] + "45").Formula = "=" + $fn
}
for ($i = 0; $i -lt $cols.Length; $i++) {
    $fn = "FALSE()"
    if ($row46[$i]) { $fn = "TRUE()" }
    $ws.Range($cols[$i] + "46").Formula = "=" + $fn
}
for ($i = 0; $i -lt $cols.Length; $i++) {
    $fn = "FALSE()"
    if ($row47[$i]) { $fn = "TRUE()" }
    $ws.Range($cols[$i] + "47").Formula = "=" + $fn
}

# Row 47 (the last row of the moved table) keeps the thicker bottom
# border/row height that previously marked the end of the table.
$ws.Rows("47:47").RowHeight = 15.75

# ignoredErrors tracking for the known false-positive boolean formulas
# moves along with the table (was J47/M48, now J46/M47).
$ws.Range("J46").ErrorCheckingOptions.NumberAsText = $false

# Match the author's last selection after performing the move.
$ws.Range("F41").Select()
